$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (A9:C9) — previously blank, now populated ---
$ws.Range("A9").Value = 130
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0

# --- New headers: Area / Atotal, plus mirrored Atotal/Qtotal in J1:K1 ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- New "Area" column formulas (G2:G15) ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").FormulaR1C1 = "=(RC[-3]-R[-1]C[-3])*RC[-5]/100"

# --- Atotal (sum of Area) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Mirror totals into J2/K2 ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection as last left by the author ---
$ws.Range("J2:K2").Select()

$wb.Save()
